# Update report header text (volume number and week date range)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Numeric value updates across the weekly crime-stat table (rows 15-31) ---
$ws.Range("N15").Value = -56.25
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = -60
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 105
$ws.Range("K16").Value = -10.476190476190
$ws.Range("L16").Value = -39.743589743589
$ws.Range("M16").Value = -40.127388535031
$ws.Range("N16").Value = -87.808041504539
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 155
$ws.Range("J17").Value = 165
$ws.Range("K17").Value = -6.060606060606
$ws.Range("L17").Value = -15.760869565217
$ws.Range("M17").Value = 23.015873015873
$ws.Range("N17").Value = -66.594827586206
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 129
$ws.Range("J18").Value = 129
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -39.150943396226
$ws.Range("M18").Value = -27.932960893854
$ws.Range("N18").Value = -82.614555256064
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -13.698630136986
$ws.Range("I19").Value = 639
$ws.Range("J19").Value = 609
$ws.Range("K19").Value = 4.926108374384
$ws.Range("L19").Value = -21.498771498771
$ws.Range("M19").Value = 4.071661237785
$ws.Range("N19").Value = -50.465116279069
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -11.428571428571
$ws.Range("L20").Value = -16.216216216216
$ws.Range("M20").Value = -24.390243902439
$ws.Range("N20").Value = -92.986425339366
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -23.931623931623
$ws.Range("I21").Value = 1063
$ws.Range("J21").Value = 1060
$ws.Range("K21").Value = 0.283018867924
$ws.Range("L21").Value = -25.193525686136
$ws.Range("M21").Value = -6.012378426171
$ws.Range("N21").Value = -71.683537559936
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -58.333333333333
$ws.Range("I23").Value = 106
$ws.Range("J23").Value = 74
$ws.Range("K23").Value = 43.243243243243
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 10.416666666666
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 13.793103448275
$ws.Range("F24").Value = 169
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = 19.014084507042
$ws.Range("I24").Value = 1283
$ws.Range("J24").Value = 1220
$ws.Range("K24").Value = 5.163934426229
$ws.Range("L24").Value = 6.208609271523
$ws.Range("M24").Value = -7.028985507246
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 5.882352941176
$ws.Range("F25").Value = 104
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 31.645569620253
$ws.Range("I25").Value = 654
$ws.Range("J25").Value = 707
$ws.Range("K25").Value = -7.496463932107
$ws.Range("L25").Value = -4.803493449781
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 359
$ws.Range("J26").Value = 345
$ws.Range("K26").Value = 4.057971014492
$ws.Range("L26").Value = -3.494623655913
$ws.Range("M26").Value = -5.774278215223
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("N29").Value = -85.714285714285
$ws.Range("N30").Value = -91.304347826087

# --- Cells that change between numeric and placeholder text ("0" / "***.*") ---
# C16: placeholder "0" -> numeric 1 (style like D16)
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1

# C17: placeholder "0" -> numeric 3 (style like D17)
$ws.Range("D17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 3

# C22: numeric 1 -> placeholder "0" (style like C14)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# D22: numeric 2 -> placeholder "0" (style like D14)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# E22: numeric -50 -> placeholder "***.*" (style like E14)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# C28: numeric 2 -> placeholder "0" (style like D28)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# G31: numeric 1 -> placeholder "0" (style like F31)
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("F31").Copy()
$ws.Range("G31").PasteSpecial(-4122)

# H31: numeric -100 -> placeholder "***.*" (style like E31)
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("E31").Copy()
$ws.Range("H31").PasteSpecial(-4122)

Write-Host "Edit complete"